$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.472.70'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '3.150.72'
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.49'
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.14'
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '3.149.15'
$ws.Range("E8").Value = '  +0.72%  '

$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.94'
$ws.Range("E11").Value = '  +3.87%  '

$ws.Range("E12").Value = '  -1.28%  '

$ws.Range("E13").Value = '  -2.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.47'
$ws.Range("E14").Value = '  +2.72%  '

$ws.Range("D15").Value = '3.674.16'
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.25'
$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("D18").Value = '64.206.37'
$ws.Range("E18").Value = '  +0.81%  '

$ws.Range("D19").Value = '3.150.70'
$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.40'
$ws.Range("E20").Value = '  +0.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.51'
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.55'
$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("E24").Value = '  -1.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +5.97%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.67'
$ws.Range("E26").Value = '  -1.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.66'
$ws.Range("E28").Value = '  +6.55%  '

$ws.Range("E29").Value = '  +2.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.49'
$ws.Range("E30").Value = '  +8.34%  '

$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.49'
$ws.Range("E33").Value = '  +0.77%  '

$ws.Range("E34").Value = '  +0.57%  '

$ws.Range("D35").Value = '0.0₃0844'
$ws.Range("E35").Value = '  -2.77%  '

$ws.Range("E36").Value = '  +1.43%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("E37").Value = '  -1.56%  '

$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.22'
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.28'
$ws.Range("E39").Value = '  -2.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.84'
$ws.Range("E40").Value = '  +1.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '458.02'
$ws.Range("E41").Value = '  +1.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.29'
$ws.Range("E42").Value = '  +5.77%  '

$ws.Range("E43").Value = '  +6.65%  '

$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D45").Value = '2.940.03'
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.96'
$ws.Range("E46").Value = '  +10.78%  '

$ws.Range("E47").Value = '  -2.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.38'
$ws.Range("E48").Value = '  +2.82%  '

$ws.Range("E50").Value = '  +2.51%  '

$ws.Range("E51").Value = '  -0.73%  '
